$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table (years 2014-2022 in columns B:J) is extended one more year to the
# right: add a "2023" column (K), copying the look of the previous year
# column (J) - same font/fill/border/number format - then fill in the new
# year header and the three data rows.
$ws.Range("J3:J6").Copy() | Out-Null
$ws.Range("K3:K6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New "2023" column values
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1392.1   # Average monthly remuneration
$ws.Range("K5").Value = 871.3    # Women
$ws.Range("K6").Value = 1512.5   # Men
